# Auto-generated edit script: update currentAveragePrice / Leve profit columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2726.3
$ws.Range("I64").Value = 2722.8572
$ws.Range("J64").Value = 2734.3333
$ws.Range("K64").Value = 2722.8572
$ws.Range("L64").Value = 2734.3333
$ws.Range("M64").Value = -2474.8572
$ws.Range("N64").Value = -3230.3333
$ws.Range("H67").Value = 2726.3
$ws.Range("I67").Value = 2722.8572
$ws.Range("J67").Value = 2734.3333
$ws.Range("K67").Value = 2722.8572
$ws.Range("L67").Value = 2734.3333
$ws.Range("M67").Value = -1864.8572
$ws.Range("N67").Value = -4450.3333
$ws.Range("H80").Value = 666.64703
$ws.Range("I80").Value = 1033.3334
$ws.Range("J80").Value = 588.0714
$ws.Range("K80").Value = 3100.0002
$ws.Range("L80").Value = 1764.2142
$ws.Range("M80").Value = -2102.0002
$ws.Range("N80").Value = -3760.2142
$ws.Range("H83").Value = 666.64703
$ws.Range("I83").Value = 1033.3334
$ws.Range("J83").Value = 588.0714
$ws.Range("K83").Value = 9300.000599999999
$ws.Range("L83").Value = 5292.6426
$ws.Range("M83").Value = -4308.000599999999
$ws.Range("N83").Value = -15276.6426
$ws.Range("H113").Value = 697013.4399999999
$ws.Range("I113").Value = 1012365
$ws.Range("J113").Value = 3240
$ws.Range("K113").Value = 1012365
$ws.Range("L113").Value = 3240
$ws.Range("M113").Value = -1009111
$ws.Range("N113").Value = -9748
$ws.Range("H116").Value = 2984444.8
$ws.Range("J116").Value = 9515
$ws.Range("L116").Value = 9515
$ws.Range("N116").Value = -16399
$ws.Range("H132").Value = 5684214.5
$ws.Range("I132").Value = 7577288.5
$ws.Range("K132").Value = 22731865.5
$ws.Range("M132").Value = -22729335.5
$ws.Range("H137").Value = 1472.2972
$ws.Range("I137").Value = 1160.3871
$ws.Range("J137").Value = 3083.8333
$ws.Range("K137").Value = 3481.1613
$ws.Range("L137").Value = 9251.499899999999
$ws.Range("M137").Value = -931.1612999999998
$ws.Range("N137").Value = -14351.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2000.4166
$ws.Range("I132").Value = 1723.9
$ws.Range("J132").Value = 2197.9285
$ws.Range("K132").Value = 5171.700000000001
$ws.Range("L132").Value = 6593.7855
$ws.Range("M132").Value = -2641.700000000001
$ws.Range("N132").Value = -11653.7855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 26500
$ws.Range("J6").Value = 26500
$ws.Range("L6").Value = 26500
$ws.Range("N6").Value = -26726
$ws.Range("H38").Value = 32450
$ws.Range("I38").Value = 20000
$ws.Range("J38").Value = 36600
$ws.Range("K38").Value = 20000
$ws.Range("L38").Value = 36600
$ws.Range("M38").Value = -19584
$ws.Range("N38").Value = -37432
$ws.Range("H51").Value = 76000
$ws.Range("J51").Value = 76000
$ws.Range("L51").Value = 76000
$ws.Range("N51").Value = -76982
$ws.Range("H55").Value = 78944.25
$ws.Range("J55").Value = 78944.25
$ws.Range("L55").Value = 78944.25
$ws.Range("N55").Value = -79490.25
$ws.Range("H64").Value = 512.125
$ws.Range("I64").Value = 785.6667
$ws.Range("J64").Value = 348
$ws.Range("K64").Value = 785.6667
$ws.Range("L64").Value = 348
$ws.Range("M64").Value = -560.6667
$ws.Range("N64").Value = -798
$ws.Range("H67").Value = 512.125
$ws.Range("I67").Value = 785.6667
$ws.Range("J67").Value = 348
$ws.Range("K67").Value = 785.6667
$ws.Range("L67").Value = 348
$ws.Range("M67").Value = -5.666699999999992
$ws.Range("N67").Value = -1908
$ws.Range("H82").Value = 16737.5
$ws.Range("I82").Value = 5300
$ws.Range("J82").Value = 18371.428
$ws.Range("K82").Value = 5300
$ws.Range("L82").Value = 18371.428
$ws.Range("M82").Value = -4917
$ws.Range("N82").Value = -19137.428
$ws.Range("H85").Value = 16737.5
$ws.Range("I85").Value = 5300
$ws.Range("J85").Value = 18371.428
$ws.Range("K85").Value = 5300
$ws.Range("L85").Value = 18371.428
$ws.Range("M85").Value = -3974
$ws.Range("N85").Value = -21023.428
$ws.Range("H86").Value = 1805.6
$ws.Range("I86").Value = 1982.6666
$ws.Range("J86").Value = 1540
$ws.Range("K86").Value = 1982.6666
$ws.Range("L86").Value = 1540
$ws.Range("M86").Value = -859.6666
$ws.Range("N86").Value = -3786
$ws.Range("H89").Value = 1805.6
$ws.Range("I89").Value = 1982.6666
$ws.Range("J89").Value = 1540
$ws.Range("K89").Value = 9913.333000000001
$ws.Range("L89").Value = 7700
$ws.Range("M89").Value = -4297.333000000001
$ws.Range("N89").Value = -18932
$ws.Range("H107").Value = 1275
$ws.Range("I107").Value = 1275
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1275
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 645
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16669942
$ws.Range("I31").Value = 31252478
$ws.Range("J31").Value = 4186.2144
$ws.Range("K31").Value = 31252478
$ws.Range("L31").Value = 4186.2144
$ws.Range("M31").Value = -31252183
$ws.Range("N31").Value = -4776.2144
$ws.Range("H34").Value = 16669942
$ws.Range("I34").Value = 31252478
$ws.Range("J34").Value = 4186.2144
$ws.Range("K34").Value = 31252478
$ws.Range("L34").Value = 4186.2144
$ws.Range("M34").Value = -31252276
$ws.Range("N34").Value = -4590.2144
$ws.Range("H38").Value = 23361.334
$ws.Range("J38").Value = 23361.334
$ws.Range("L38").Value = 23361.334
$ws.Range("N38").Value = -24115.334
$ws.Range("H46").Value = 23361.334
$ws.Range("J46").Value = 23361.334
$ws.Range("L46").Value = 23361.334
$ws.Range("N46").Value = -23783.334
$ws.Range("H80").Value = 8595
$ws.Range("J80").Value = 8595
$ws.Range("L80").Value = 8595
$ws.Range("N80").Value = -10841
$ws.Range("H83").Value = 8595
$ws.Range("J83").Value = 8595
$ws.Range("L83").Value = 25785
$ws.Range("N83").Value = -37017
$ws.Range("H107").Value = 1352.6875
$ws.Range("I107").Value = 664.8461
$ws.Range("J107").Value = 4333.3335
$ws.Range("K107").Value = 664.8461
$ws.Range("L107").Value = 4333.3335
$ws.Range("M107").Value = 1255.1539
$ws.Range("N107").Value = -8173.3335
$ws.Range("H132").Value = 3704.1538
$ws.Range("I132").Value = 2560
$ws.Range("J132").Value = 4419.25
$ws.Range("K132").Value = 7680
$ws.Range("L132").Value = 13257.75
$ws.Range("M132").Value = -5150
$ws.Range("N132").Value = -18317.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3333.2222
$ws.Range("J88").Value = 3333.2222
$ws.Range("L88").Value = 9999.6666
$ws.Range("N88").Value = -10855.6666
$ws.Range("H91").Value = 3333.2222
$ws.Range("J91").Value = 3333.2222
$ws.Range("L91").Value = 9999.6666
$ws.Range("N91").Value = -12963.6666
$ws.Range("H113").Value = 396.15384
$ws.Range("I113").Value = 356.89655
$ws.Range("J113").Value = 510
$ws.Range("K113").Value = 1070.68965
$ws.Range("L113").Value = 1530
$ws.Range("M113").Value = 1099.31035
$ws.Range("N113").Value = -5870
$ws.Range("H131").Value = 838.2545
$ws.Range("I131").Value = 382.9
$ws.Range("J131").Value = 939.44446
$ws.Range("K131").Value = 1148.7
$ws.Range("L131").Value = 2818.33338
$ws.Range("M131").Value = 3891.3
$ws.Range("N131").Value = -12898.33338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6134.5386
$ws.Range("I70").Value = 5033.222
$ws.Range("J70").Value = 8612.5
$ws.Range("K70").Value = 5033.222
$ws.Range("L70").Value = 8612.5
$ws.Range("M70").Value = -4763.222
$ws.Range("N70").Value = -9152.5
$ws.Range("H73").Value = 6134.5386
$ws.Range("I73").Value = 5033.222
$ws.Range("J73").Value = 8612.5
$ws.Range("K73").Value = 5033.222
$ws.Range("L73").Value = 8612.5
$ws.Range("M73").Value = -4097.222
$ws.Range("N73").Value = -10484.5
$ws.Range("H107").Value = 1807.5555
$ws.Range("I107").Value = 2167.3333
$ws.Range("J107").Value = 1627.6666
$ws.Range("K107").Value = 2167.3333
$ws.Range("L107").Value = 1627.6666
$ws.Range("M107").Value = -247.3332999999998
$ws.Range("N107").Value = -5467.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 18160
$ws.Range("J94").Value = 18160
$ws.Range("L94").Value = 18160
$ws.Range("N94").Value = -19512

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 33659.8
$ws.Range("J93").Value = 33659.8
$ws.Range("L93").Value = 33659.8
$ws.Range("N93").Value = -38651.8
